$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.006434674933132
$ws.Range("C2").Value = 1.564547541784123
$ws.Range("D2").Value = 0.5283182698003515
$ws.Range("E2").Value = 0.1640833435699314
$ws.Range("G2").Value = 0.002666014655108295
$ws.Range("I2").Value = 5.445869992638507
$ws.Range("J2").Value = 0.05312248822832899
$ws.Range("M2").Value = 1.029616708816015

$ws.Range("B3").Value = 1.97642842872915
$ws.Range("C3").Value = 1.503609022336718
$ws.Range("D3").Value = 0.5251348404127754
$ws.Range("E3").Value = 0.16339048030326
$ws.Range("G3").Value = 0.002674642672323952
$ws.Range("I3").Value = 5.125513018071189
$ws.Range("J3").Value = 0.05295856339539817
$ws.Range("M3").Value = 1.012928523712695

$ws.Range("B4").Value = 1.960548681672378
$ws.Range("C4").Value = 1.467254053891168
$ws.Range("D4").Value = 0.5234232275861643
$ws.Range("E4").Value = 0.1630230955894163
$ws.Range("G4").Value = 0.002680200693049759
$ws.Range("I4").Value = 4.92886761268241
$ws.Range("J4").Value = 0.05286446046315341
$ws.Range("M4").Value = 1.003558908797949

$ws.Range("B5").Value = 1.954712125084598
$ws.Range("C5").Value = 1.452702397485382
$ws.Range("D5").Value = 0.5227865767805042
$ws.Range("E5").Value = 0.1628879213927732
$ws.Range("G5").Value = 0.002682531410546535
$ws.Range("I5").Value = 4.848730752452809
$ws.Range("J5").Value = 0.05282774782088495
$ws.Range("M5").Value = 0.9999596210822546

$ws.Range("B6").Value = 1.95378114661645
$ws.Range("C6").Value = 1.450301903172317
$ws.Range("D6").Value = 0.5226845287269697
$ws.Range("E6").Value = 0.1628663522972822
$ws.Range("G6").Value = 0.002682922406255608
$ws.Range("I6").Value = 4.835423485064439
$ws.Range("J6").Value = 0.05282175014129109
$ws.Range("M6").Value = 0.9993751369831756

$ws.Range("B7").Value = 1.960467405225472
$ws.Range("C7").Value = 1.46705674386584
$ws.Range("D7").Value = 0.5234143954726278
$ws.Range("E7").Value = 0.1630212137905573
$ws.Range("G7").Value = 0.002680231859431702
$ws.Range("I7").Value = 4.927786889246477
$ws.Range("J7").Value = 0.05286395873822336
$ws.Range("M7").Value = 1.003509483305436

$ws.Range("B8").Value = 1.995557284014637
$ws.Range("C8").Value = 1.543313381779342
$ws.Range("D8").Value = 0.527170018786606
$ws.Range("E8").Value = 0.1638323643815873
$ws.Range("G8").Value = 0.002668935736343267
$ws.Range("I8").Value = 5.33538962545768
$ws.Range("J8").Value = 0.05306460099815702
$ws.Range("M8").Value = 1.023679603201565

$ws.Range("B9").Value = 2.084804451194771
$ws.Range("C9").Value = 1.701452044146151
$ws.Range("D9").Value = 0.5364763981491478
$ws.Range("E9").Value = 0.1658863317710484
$ws.Range("G9").Value = 0.002648836170086496
$ws.Range("I9").Value = 6.135904502400365
$ws.Range("J9").Value = 0.05351056879250038
$ws.Range("M9").Value = 1.070270671636933

$ws.Range("B10").Value = 2.163199038461983
$ws.Range("C10").Value = 1.823142547529073
$ws.Range("D10").Value = 0.5445178637922368
$ws.Range("E10").Value = 0.1676821678725489
$ws.Range("G10").Value = 0.002635300397056122
$ws.Range("I10").Value = 6.725908999222895
$ws.Range("J10").Value = 0.05387106763681437
$ws.Range("M10").Value = 1.108909239930284

$ws.Range("B11").Value = 2.201730880565663
$ws.Range("C11").Value = 1.879757996801743
$ws.Range("D11").Value = 0.5484423024497858
$ws.Range("E11").Value = 0.1685624193933641
$ws.Range("G11").Value = 0.00262940575539413
$ws.Range("I11").Value = 6.994975274094315
$ws.Range("J11").Value = 0.05404238978973908
$ws.Range("M11").Value = 1.127471011547328

$ws.Range("B12").Value = 2.21674141798195
$ws.Range("C12").Value = 1.901382434901677
$ws.Range("D12").Value = 0.5499670610532235
$ws.Range("E12").Value = 0.168904932632703
$ws.Range("G12").Value = 0.00262721107188169
$ws.Range("I12").Value = 7.096980068301832
$ws.Range("J12").Value = 0.05410833403263382
$ws.Range("M12").Value = 1.134643679677154

$ws.Range("B13").Value = 2.213489869701334
$ws.Range("C13").Value = 1.896716906252777
$ws.Range("D13").Value = 0.5496369517617836
$ws.Range("E13").Value = 0.1688307567011265
$ws.Range("G13").Value = 0.002627682073953027
$ws.Range("I13").Value = 7.075006080857747
$ws.Range("J13").Value = 0.0540940840167643
$ws.Range("M13").Value = 1.133092490672695

$ws.Range("B14").Value = 2.202957362647055
$ws.Range("C14").Value = 1.881533310332486
$ws.Range("D14").Value = 0.5485669683643835
$ws.Range("E14").Value = 0.168590413713293
$ws.Range("G14").Value = 0.002629224447764134
$ws.Range("I14").Value = 7.003364850887067
$ws.Range("J14").Value = 0.05404779356942058
$ws.Range("M14").Value = 1.128058218450875

$ws.Range("B15").Value = 2.196560714300063
$ws.Range("C15").Value = 1.872257200014246
$ws.Range("D15").Value = 0.5479166180917332
$ws.Range("E15").Value = 0.1684443946576017
$ws.Range("G15").Value = 0.002630174069577577
$ws.Range("I15").Value = 6.95949810286902
$ws.Range("J15").Value = 0.05401957888438957
$ws.Range("M15").Value = 1.124993362599298

$ws.Range("B16").Value = 2.16073920616958
$ws.Range("C16").Value = 1.81946828896173
$ws.Range("D16").Value = 0.5442667823149918
$ws.Range("E16").Value = 0.1676259219066658
$ws.Range("G16").Value = 0.002635690888688484
$ws.Range("I16").Value = 6.708339974618525
$ws.Range("J16").Value = 0.05386002003058366
$ws.Range("M16").Value = 1.10771618108086

$ws.Range("B17").Value = 2.139503578649567
$ws.Range("C17").Value = 1.787409336778296
$ws.Range("D17").Value = 0.5420961930534816
$ws.Range("E17").Value = 0.1671400855936582
$ws.Range("G17").Value = 0.002639142382585143
$ws.Range("I17").Value = 6.554447548463543
$ws.Range("J17").Value = 0.05376402312385586
$ws.Range("M17").Value = 1.0973709172479

$ws.Range("B18").Value = 2.127559107867
$ws.Range("C18").Value = 1.769087989614547
$ws.Range("D18").Value = 0.5408727739066705
$ws.Range("E18").Value = 0.1668666003200272
$ws.Range("G18").Value = 0.00264115234923066
$ws.Range("I18").Value = 6.46599463108484
$ws.Range("J18").Value = 0.05370949724168739
$ws.Range("M18").Value = 1.091513167895201

$ws.Range("B19").Value = 2.123561036109322
$ws.Range("C19").Value = 1.762904846937545
$ws.Range("D19").Value = 0.5404628362018684
$ws.Range("E19").Value = 0.1667750233960064
$ws.Range("G19").Value = 0.002641837151570639
$ws.Range("I19").Value = 6.436056063251982
$ws.Range("J19").Value = 0.05369115362267252
$ws.Range("M19").Value = 1.089545671202956

$ws.Range("B20").Value = 2.141736186674962
$ws.Range("C20").Value = 1.790809810271867
$ws.Range("D20").Value = 0.5423246612157584
$ws.Range("E20").Value = 0.1671911868440539
$ws.Range("G20").Value = 0.002638772404951984
$ws.Range("I20").Value = 6.570823114623352
$ws.Range("J20").Value = 0.05377417074770463
$ws.Range("M20").Value = 1.098462592320274

$ws.Range("B21").Value = 2.20603957932326
$ws.Range("C21").Value = 1.885988032197304
$ws.Range("D21").Value = 0.5488801962953858
$ws.Range("E21").Value = 0.1686607584800726
$ws.Range("G21").Value = 0.002628770400298902
$ws.Range("I21").Value = 7.024404336899465
$ws.Range("J21").Value = 0.05406136108714321
$ws.Range("M21").Value = 1.129532988386671

$ws.Range("B22").Value = 2.250513215180433
$ws.Range("C22").Value = 1.949275330162607
$ws.Range("D22").Value = 0.5533901280957707
$ws.Range("E22").Value = 0.1696747604997917
$ws.Range("G22").Value = 0.002622451886115791
$ws.Range("I22").Value = 7.321526627196192
$ws.Range("J22").Value = 0.05425529016921615
$ws.Range("M22").Value = 1.150678072990729

$ws.Range("B23").Value = 2.226550566081471
$ws.Range("C23").Value = 1.915397099466304
$ws.Range("D23").Value = 0.5509623356999498
$ws.Range("E23").Value = 0.1691286425424465
$ws.Range("G23").Value = 0.00262580431684918
$ws.Range("I23").Value = 7.162878239131032
$ws.Range("J23").Value = 0.05415121135757772
$ws.Range("M23").Value = 1.139315081135052

$ws.Range("B24").Value = 2.140726002246936
$ws.Range("C24").Value = 1.789272114631785
$ws.Range("D24").Value = 0.5422212945307763
$ws.Range("E24").Value = 0.1671680658200607
$ws.Range("G24").Value = 0.002638939591570689
$ws.Range("I24").Value = 6.563419654405408
$ws.Range("J24").Value = 0.05376958094007733
$ws.Range("M24").Value = 1.097968766219978

$ws.Range("B25").Value = 2.058435765376998
$ws.Range("C25").Value = 1.657723899321127
$ws.Range("D25").Value = 0.5337488485819648
$ws.Range("E25").Value = 0.1652806400934423
$ws.Range("G25").Value = 0.002654055968716018
$ws.Range("I25").Value = 5.919099307260097
$ws.Range("J25").Value = 0.05338421623342882
$ws.Range("M25").Value = 1.056901188818713

